# New classifications in Test Dataset
#
# The "Teste" sheet has two columns: A = classification (0/1/2), B = text.
# Rows 76-154 previously had no real classification (row 76 held a lone
# space placeholder, rows 77-154 had no value at all in column A). This
# fills in the newly-assigned numeric classifications for that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teste")

$values = @(0,2,2,1,2,2,0,2,2,1,1,0,1,0,2,0,2,0,1,2,2,2,1,2,2,0,0,1,1,2,2,2,2,2,2,0,2,1,2,1,0,2,1,1,2,1,2,2,0,1,0,2,2,2,2,1,2,0,1,0,0,2,2,2,2,0,1,1,1,2,0,0,2,0,1,2,1,2,2)

$startRow = 76
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Restore the author's final cursor position on the sheet (matches the
# sheetView selection recorded after entering the new data).
$ws.Range("A155").Select()
